$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2023-09-25) needs to be inserted above the
# existing row 92, pushing the historical rows (92-97) down to (93-98).
$ws.Rows("92:92").Insert()

# Populate the newly inserted row 92 with the new week's data.
$ws.Cells.Item(92, 1).Value = 3
$ws.Cells.Item(92, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(92, 3).Value = "Coquimbo"
$ws.Cells.Item(92, 4).Value = 45194
$ws.Cells.Item(92, 5).Value = 5
$ws.Cells.Item(92, 6).Value = 100112022
$ws.Cells.Item(92, 7).Value = "Arveja Verde"
$ws.Cells.Item(92, 8).Value = "Perfection"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 40
$ws.Cells.Item(92, 11).Value = 32000
$ws.Cells.Item(92, 12).Value = 32000
$ws.Cells.Item(92, 13).Value = 32000
$ws.Cells.Item(92, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(92, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(92, 16).Value = 1280
$ws.Cells.Item(92, 17).Value = 25
$ws.Cells.Item(92, 18).Value = "Hortaliza"
